$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 41.22222
$ws.Range("I6").Value = 35.4
$ws.Range("K6").Value = 106.2
$ws.Range("M6").Value = 5.800000000000011

$ws.Range("H19").Value = 291.18182
$ws.Range("I19").Value = 336.33334
$ws.Range("J19").Value = 88
$ws.Range("K19").Value = 336.33334
$ws.Range("L19").Value = 88
$ws.Range("M19").Value = -161.33334
$ws.Range("N19").Value = -438

$ws.Range("H112").Value = 1999.9
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1999.9
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 5999.700000000001
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -8215.700000000001

$ws.Range("H118").Value = 162.8
$ws.Range("I118").Value = 162.8
$ws.Range("K118").Value = 488.4
$ws.Range("M118").Value = 1168.6

$ws.Range("H135").Value = 501.3
$ws.Range("I135").Value = 420.54544
$ws.Range("J135").Value = 600
$ws.Range("K135").Value = 3784.90896
$ws.Range("L135").Value = 5400
$ws.Range("M135").Value = -1249.90896
$ws.Range("N135").Value = -10470

$ws.Range("H138").Value = 3052.7273
$ws.Range("I138").Value = 1500
$ws.Range("J138").Value = 3208
$ws.Range("K138").Value = 4500
$ws.Range("L138").Value = 9624
$ws.Range("M138").Value = 640
$ws.Range("N138").Value = -19904

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 184.4762
$ws.Range("I5").Value = 185.5
$ws.Range("K5").Value = 185.5
$ws.Range("M5").Value = -73.5

$ws.Range("H32").Value = 3141.7058
$ws.Range("I32").Value = 3141.7058
$ws.Range("K32").Value = 3141.7058
$ws.Range("M32").Value = -2854.7058

$ws.Range("H61").Value = 1536.75
$ws.Range("I61").Value = 1091.5
$ws.Range("J61").Value = 1982
$ws.Range("K61").Value = 1091.5
$ws.Range("L61").Value = 1982
$ws.Range("M61").Value = -879.5
$ws.Range("N61").Value = -2406

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H97").Value = 393.63635
$ws.Range("I97").Value = 393.63635
$ws.Range("K97").Value = 393.63635
$ws.Range("M97").Value = 102.36365

$ws.Range("H107").Value = 70000
$ws.Range("J107").Value = 70000
$ws.Range("L107").Value = 70000
$ws.Range("N107").Value = -77680

$ws.Range("H132").Value = 3169.6667
$ws.Range("J132").Value = 3985
$ws.Range("L132").Value = 11955
$ws.Range("N132").Value = -17015

$ws.Range("H136").Value = 1536.75
$ws.Range("I136").Value = 1091.5
$ws.Range("J136").Value = 1982
$ws.Range("K136").Value = 3274.5
$ws.Range("L136").Value = 5946
$ws.Range("M136").Value = -724.5
$ws.Range("N136").Value = -11046

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 184.4762
$ws.Range("I4").Value = 185.5
$ws.Range("K4").Value = 185.5
$ws.Range("M4").Value = -70.5

$ws.Range("H86").Value = 2244.2727
$ws.Range("I86").Value = 855.2857
$ws.Range("J86").Value = 4675
$ws.Range("K86").Value = 855.2857
$ws.Range("L86").Value = 4675
$ws.Range("M86").Value = 267.7143
$ws.Range("N86").Value = -6921

$ws.Range("H89").Value = 2244.2727
$ws.Range("I89").Value = 855.2857
$ws.Range("J89").Value = 4675
$ws.Range("K89").Value = 4276.4285
$ws.Range("L89").Value = 23375
$ws.Range("M89").Value = 1339.5715
$ws.Range("N89").Value = -34607

$ws.Range("H134").Value = 3242.5
$ws.Range("J134").Value = 3485
$ws.Range("L134").Value = 10455
$ws.Range("N134").Value = -15525

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 230000
$ws.Range("J9").Value = 230000
$ws.Range("L9").Value = 230000
$ws.Range("N9").Value = -230336

$ws.Range("H22").Value = 732.8333
$ws.Range("I22").Value = 933.7692
$ws.Range("J22").Value = 495.36365
$ws.Range("K22").Value = 933.7692
$ws.Range("L22").Value = 495.36365
$ws.Range("M22").Value = -583.7692
$ws.Range("N22").Value = -1195.36365

$ws.Range("H58").Value = 666.3333
$ws.Range("I58").Value = 499.5
$ws.Range("K58").Value = 499.5
$ws.Range("M58").Value = -296.5

$ws.Range("H107").Value = 721
$ws.Range("I107").Value = 705.4
$ws.Range("K107").Value = 705.4
$ws.Range("M107").Value = 1214.6

$ws.Range("H136").Value = 666.3333
$ws.Range("I136").Value = 499.5
$ws.Range("K136").Value = 1498.5
$ws.Range("M136").Value = 1051.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4986.5
$ws.Range("I80").Value = 4986.5
$ws.Range("K80").Value = 14959.5
$ws.Range("M80").Value = -14023.5

$ws.Range("H83").Value = 4986.5
$ws.Range("I83").Value = 4986.5
$ws.Range("K83").Value = 44878.5
$ws.Range("M83").Value = -40198.5

$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 60000
$ws.Range("N88").Value = -60856

$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 60000
$ws.Range("N91").Value = -62964

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 12000
$ws.Range("J49").Value = 12000
$ws.Range("L49").Value = 12000
$ws.Range("N49").Value = -12368

$ws.Range("H132").Value = 2675.3333
$ws.Range("I132").Value = 2012
$ws.Range("J132").Value = 3007
$ws.Range("K132").Value = 6036
$ws.Range("L132").Value = 9021
$ws.Range("M132").Value = -3506
$ws.Range("N132").Value = -14081

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 914.6667
$ws.Range("I55").Value = 704.3333
$ws.Range("K55").Value = 704.3333
$ws.Range("M55").Value = -531.3333

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H136").Value = 3400
$ws.Range("J136").Value = 2800
$ws.Range("L136").Value = 8400
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
